$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: lowercase / snake_case renames ---
$ws.Range("A1").Value = "fecha"
$ws.Range("B1").Value = "proveedor"
$ws.Range("C1").Value = "producto"
$ws.Range("D1").Value = "cantidad"
$ws.Range("E1").Value = "precio_unitario"
$ws.Range("F1").Value = "precio_total"

# --- Column A: convert date-serial cells to plain "YYYY-MM-DD" text, dropping the old number-format style ---
$dates = @{
    2  = "2025-05-10"
    3  = "2025-05-12"
    4  = "2025-05-14"
    5  = "2025-05-17"
    6  = "2025-05-12"
    7  = "2025-05-14"
    8  = "2025-05-15"
    9  = "2025-05-18"
    10 = "2025-05-14"
    11 = "2025-05-16"
    12 = "2025-05-17"
    13 = "2025-05-20"
    14 = "2025-05-13"
    15 = "2025-05-15"
    16 = "2025-05-17"
    17 = "2025-05-18"
}

foreach ($r in $dates.Keys) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.ClearContents()
    $cell.ClearFormats()
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$r]
    $cell.Style = "Normal"
}

# --- Column B: trim trailing spaces from provider names ---
$ws.Range("B2").Value = "LogiMax"
$ws.Range("B3").Value = "LogiMax"
$ws.Range("B4").Value = "LogiMax"
$ws.Range("B5").Value = "LogiMax"
$ws.Range("B10").Value = "NovaIndustrias"
$ws.Range("B11").Value = "NovaIndustrias"
$ws.Range("B12").Value = "NovaIndustrias"
$ws.Range("B13").Value = "NovaIndustrias"
